# Adds the "Agendamento", "Cadastro Tutor", "Cadastro animal" and
# "Estoque" scenarios (rows 6-15) to the test-script sheet, reusing the
# existing header/body cell formatting already present on the sheet
# (border + left/center/wrap for the bulk of the table, and the
# border-less wrap-only look used on the final "Estoque" row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use row 2's formatting (border on all sides, left/top aligned, wrapped
# text) as the template for the new data rows, the same way the rest of
# the sheet's rows already look.

$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A6:E6").PasteSpecial(-4122) | Out-Null
$ws.Range("A6").Value = 'Agendamento'
$ws.Range("B6").Value = 'Funcionário realizando um agendamento'
$ws.Range("C6").Value = '1. Funcionário entra na aba de "Agendamentos"' + "`n" + '2. Clicar no dia 21 de stembro de 2022 e no horário das 18:30, que é o momento escolhido para o agendamento' + "`n" + '3. Clicar no botão "Novo agendamento"' + "`n" + '4. Preencher os dados do veterinário Luís, do tutor João, do animal cachorro e a consulta de rotina' + "`n" + '5. Clicar no botão de salvar'
$ws.Range("D6").Value = '1. Funcionário entra na aba de "Agendamentos"' + "`n" + '2. Clicar no dia e no horário que deverá ser feito o agendamento' + "`n" + '3. Clicar no botão "Novo agendamento"' + "`n" + '4. Preencher os dados do veterinário, tutor, animal e o tipo de consulta' + "`n" + '5. Clicar no botão de salvar'
$ws.Rows(6).RowHeight = 75

$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A7:E7").PasteSpecial(-4122) | Out-Null
$ws.Range("A7").Value = 'Agendamento'
$ws.Range("B7").Value = 'Desistência de marcar um agendamento'
$ws.Range("C7").Value = '1. Funcionário entra na aba de "Agendamentos"' + "`n" + '2. Clicar no dia 21 de stembro de 2022 e no horário das 18:30, que é o momento escolhido para o agendamento' + "`n" + '3. Clicar no botão "Novo agendamento"' + "`n" + '4. Preencher os dados do veterinário Luís, do tutor João, do animal cachorro e a consulta de rotina' + "`n" + '5. Há uma desistência no agendamento, então não clica no botão de salvar e sai da página'
$ws.Range("D7").Value = '1. Funcionário entra na aba de "Agendamentos"' + "`n" + '2. Clicar no dia e no horário que deverá ser feito o agendamento' + "`n" + '3. Clicar no botão "Novo agendamento"' + "`n" + '4. Preencher os dados do veterinário, tutor, animal e o tipo de consulta' + "`n" + '5. Não clica no botão de salvar e sai da página'
$ws.Rows(7).RowHeight = 75

$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A8:E8").PasteSpecial(-4122) | Out-Null
$ws.Range("A8").Value = 'Agendamento'
$ws.Range("B8").Value = 'Mudança de dia e horário de um agendamento'
$ws.Range("C8").Value = '1. Funcionário entra na aba de "Agendamentos"' + "`n" + '2. Clicar no dia 16 de agosto de 2022 e no horário das 15:00, que é o dia e horário anteriormente marcado para a consulta' + "`n" + '3. Clicar no botão "Editar"' + "`n" + '4. Altera os dados do veterinário Luís, do tutor João, do animal cachorro e a consulta de rotina' + "`n" + '5. Clica no botão salvar'
$ws.Range("D8").Value = '1. Funcionário entra na aba de "Agendamentos"' + "`n" + '2. Clicar no dia e no horário que o agendamento foi feito anteriormente' + "`n" + '3. Clicar no botão "Editar"' + "`n" + '4. Preencher os novos dados do veterinário, tutor, animal e o tipo de consulta' + "`n" + '5. Clica no botão de salvar'
$ws.Rows(8).RowHeight = 75

$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A9:E9").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Value = 'Agendamento'
$ws.Range("B9").Value = 'Cancelamento de um agendamento'
$ws.Range("C9").Value = '1. Funcionário entra na aba de "Agendamentos"' + "`n" + '2. Clicar no dia 16 de agosto de 2022 e no horário das 15:00, que é o dia e horário anteriormente marcado para a consulta' + "`n" + '3. Clicar no botão "Editar"' + "`n" + '4. Clica no botão "Cancelar"'
$ws.Range("D9").Value = '1. Funcionário entra na aba de "Agendamentos"' + "`n" + '2. Clicar no dia e no horário que o agendamento foi feito anteriormente' + "`n" + '3. Clicar no botão "Editar"' + "`n" + '4. Clica no botão de "Cancelar"'
$ws.Rows(9).RowHeight = 75

$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A10:E10").PasteSpecial(-4122) | Out-Null
$ws.Range("A10").Value = 'Agendamento'
$ws.Range("B10").Value = 'Cancelamento de um agendamento'
$ws.Range("C10").Value = '1. Funcionário entra na aba de "Agendamentos"' + "`n" + '2. Clicar no dia 16 de agosto de 2022 e no horário das 15:00, que é o dia e horário anteriormente marcado para a consulta' + "`n" + '3. Clicar no botão "Editar"' + "`n" + '4. Clica no botão "Cancelar"'
$ws.Range("D10").Value = '1. Funcionário entra na aba de "Agendamentos"' + "`n" + '2. Clicar no dia e no horário que o agendamento foi feito anteriormente' + "`n" + '3. Clicar no botão "Editar"' + "`n" + '4. Clica no botão de "Cancelar"'
$ws.Rows(10).RowHeight = 75

$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A11:E11").PasteSpecial(-4122) | Out-Null
$ws.Range("A11").Value = 'Cadastro Tutor'
$ws.Range("B11").Value = 'Cadastro de novo tutor'
$ws.Range("C11").Value = '1. Ir na aba "Cadastro tutor"' + "`n" + '2. Clicar no botão "Novo cadastro"' + "`n" + '3. Preencher os dados de nome Matheus, Endereço Alameda dos Arapanés 1.441, CPF 012.032.761-89, Telefone (11) 3726-6207, E-mail matheus.leal@hotmail.com' + "`n" + '4. Clicar em "Salvar"'
$ws.Range("D11").Value = '1. Ir na aba "Cadastro tutor"' + "`n" + '2. Clicar no botão "Novo cadastro"' + "`n" + '3. Preencher os dados de do tutor' + "`n" + '4. Clicar em "Salvar"'
$ws.Rows(11).RowHeight = 75

$ws.Range("A2:D2").Copy() | Out-Null
$ws.Range("A12:D12").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").Value = 'Cadastro Tutor'
$ws.Range("B12").Value = 'Editar cadastro de um tutor'
$ws.Range("C12").Value = '1. Ir na aba "Cadastro tutor"' + "`n" + '2. Escrever o nome do tutor João' + "`n" + '3. Clicar no botão "Pesquisar"' + "`n" + '4. Procurar dentre todos os registros de João, o tutor que se deve editar o cadastro, João Almeida Silva Pereira' + "`n" + '5. Clicar no registro, fazendo com que esse abra' + "`n" + '6. Clicar em "Editar"' + "`n" + '7. Alterar os dados necessários' + "`n" + '8. Clicar em "Salvar"'
$ws.Range("D12").Value = '1. Ir na aba "Cadastro tutor"' + "`n" + '2. Escrever o nome do tutor cujos dados deverão ser alterados' + "`n" + '3. Clicar no botão "Pesquisar"' + "`n" + '4. Procurar dentre todos os registros que apareceram na pesquisa, o tutor que se deve editar o cadastro' + "`n" + '5. Clicar no registro, fazendo com que esse abra' + "`n" + '6. Clicar em "Editar"' + "`n" + '7. Alterar os dados necessários' + "`n" + '8. Clicar em "Salvar"'
$ws.Rows(12).RowHeight = 180

$ws.Range("A2:D2").Copy() | Out-Null
$ws.Range("A13:D13").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").Value = 'Cadastro animal'
$ws.Range("B13").Value = 'Cadastro de novo animal'
$ws.Range("C13").Value = '1. Ir na aba "Cadastro animal"' + "`n" + '2. Clicar no botão "Novo cadastro"' + "`n" + '3. Preencher os dados de nome Lady, tutor Fernanda, pelagem preto, animal cachorro, sexo fêmea' + "`n" + '4. Clicar em "Salvar"'
$ws.Range("D13").Value = '1. Ir na aba "Cadastro animal"' + "`n" + '2. Clicar no botão "Novo cadastro"' + "`n" + '3. Preencher os dados do animal a ser cadastrado' + "`n" + '4. Clicar em "Salvar"'
$ws.Rows(13).RowHeight = 90

$ws.Range("A2:D2").Copy() | Out-Null
$ws.Range("A14:D14").PasteSpecial(-4122) | Out-Null
$ws.Range("A14").Value = 'Cadastro animal'
$ws.Range("B14").Value = 'Editar cadastro de um animal'
$ws.Range("C14").Value = '1. Ir na aba "Cadastro animal"' + "`n" + '2. Escrever o nome do animal Lady' + "`n" + '3. Clicar no botão "Pesquisar"' + "`n" + '4. Procurar dentre todos os registros de Lady, o animal que se deve editar o cadastro' + "`n" + '5. Clicar no registro, fazendo com que esse abra' + "`n" + '6. Clicar em "Editar"' + "`n" + '7. Alterar os dados necessários' + "`n" + '8. Clicar em "Salvar"'
$ws.Range("D14").Value = '1. Ir na aba "Cadastro animal"' + "`n" + '2. Escrever o nome do animal cujos dados deverão ser alterados' + "`n" + '3. Clicar no botão "Pesquisar"' + "`n" + '4. Procurar dentre todos os registros que apareceram na pesquisa, o animal que se deve editar o cadastro' + "`n" + '5. Clicar no registro, fazendo com que esse abra' + "`n" + '6. Clicar em "Editar"' + "`n" + '7. Alterar os dados necessários' + "`n" + '8. Clicar em "Salvar"'
$ws.Rows(14).RowHeight = 180

$ws.Range("A2:D2").Copy() | Out-Null
$ws.Range("A15:D15").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Value = 'Estoque'
$ws.Range("B15").Value = 'Verificar quantidade disponível de um medicamento no estoque'
$ws.Range("C15").Value = '1. Ir na aba "Estoque"' + "`n" + '2. Escrever no campo de busca o medicamento "Bayer Austria GmbH Herbststraße 6-10 1160 Wien AUSTRIA"' + "`n" + '3. Clicar no botão "Pesquisar"' + "`n" + '4. Encontre, dentre os resultados, aquele desejado' + "`n" + '5. Clique no registro desejado' + "`n" + '6. Abrirá um campo com maiores informações sobre o medicamento' + "`n" + '7. Clique no botão "Verificar estoque"' + "`n" + '8. Informações sobre a quantidade disponível e o histórico de entradas e saídas do medicamento aparecerá'
$ws.Range("D15").Value = '1. Ir na aba "Estoque"' + "`n" + '2. Escrever no campo de busca o nome comercial ou o nome genérico do medicamento no qual se deseja saber a quantidade desponível' + "`n" + '3. Clicar no botão "Pesquisar"' + "`n" + '4. Encontre, dentre os resultados, aquele desejado' + "`n" + '5. Clique no registro desejado' + "`n" + '6. Abrirá um campo com maiores informações sobre o medicamento escolhido' + "`n" + '7. Clique no botão "Verificar estoque"' + "`n" + '8. Informações sobre a quantidade disponível e o histórico de entradas e saídas do medicamento aparecerá'
$ws.Rows(15).RowHeight = 240

# Row 15 ("Estoque" / verificar estoque) loses its border in the
# source file and keeps only the wrap formatting; B15 still keeps the
# left/center alignment while C15:D15 drop back to default alignment.
$ws.Range("B15:D15").Borders.LineStyle = -4142
$ws.Range("B15:D15").Interior.Pattern = -4142
$ws.Range("C15:D15").HorizontalAlignment = 1
$ws.Range("C15:D15").VerticalAlignment = -4107

# Match the selection the workbook was left with when it was saved.
$ws.Range("I14").Select() | Out-Null
